$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (row 11) since the new table only spans rows 1-10
$ws.Rows(11).Delete()

# Update match data rows 2-10 with new teams and scores
$ws.Range("A2").Value = "Wizards"
$ws.Range("B2").Value = "Magic"
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = "Lakers"
$ws.Range("B3").Value = "Hawks"
$ws.Range("C3").Value = 4

$ws.Range("A4").Value = "Suns"
$ws.Range("B4").Value = "Raptors"
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = "Pistons"
$ws.Range("B5").Value = "Bulls"
$ws.Range("C5").Value = 3

$ws.Range("A6").Value = "Timberwolves"
$ws.Range("B6").Value = "Bucks"
$ws.Range("C6").Value = 4

$ws.Range("A7").Value = "76ers"
$ws.Range("B7").Value = "Pelicans"
$ws.Range("C7").Value = 3

$ws.Range("A8").Value = "Heat"
$ws.Range("B8").Value = "Nuggets"
$ws.Range("C8").Value = 2

$ws.Range("A9").Value = "Trail Blazers"
$ws.Range("B9").Value = "Warriors"
$ws.Range("C9").Value = 4

$ws.Range("A10").Value = "Jazz"
$ws.Range("B10").Value = "Kings"
$ws.Range("C10").Value = 4
